# Normalize accented Portuguese characters in the "NOME_UNIDADE" (B) and
# "NIVEL_CURSO" (C) columns to their plain-ASCII equivalents, e.g.
# "Graduação" -> "Graduacao", "Ciência" -> "Ciencia", etc.
# Only the data rows (2..325) are affected; the header row (row 1) has no
# accented characters so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$rng = $ws.Range("B2:C$lastRow")

$rng.Replace("á", "a") | Out-Null
$rng.Replace("â", "a") | Out-Null
$rng.Replace("ã", "a") | Out-Null
$rng.Replace("ç", "c") | Out-Null
$rng.Replace("é", "e") | Out-Null
$rng.Replace("ê", "e") | Out-Null
$rng.Replace("í", "i") | Out-Null
$rng.Replace("õ", "o") | Out-Null
$rng.Replace("ú", "u") | Out-Null
